# excel/breakout/chartink_screener.xlsx — "break out stock.yaml completed"
#
# The sheet named "DND 3 V 0.3" gets:
#   1. E3's bsecode value converted from text ("532832") to a real number (532832)
#   2. A new row 4 appended with the next screener snapshot for IBREALEST,
#      keeping bsecode (E4) as text "532832" (matching the original scraper
#      output format used for every other row in this sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DND 3 V 0.3")

# 1) E3: text "532832" -> numeric 532832
$ws.Range("E3").Value = 532832

# 2) Append row 4
$ws.Range("A4").Value = "20/06/2024 06:44:51"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "IBREALEST"
$ws.Range("D4").Value = "Indiabulls Real Estate Limited"

# Force E4 to stay text (bsecode "532832") rather than auto-converting to a
# number: a leading apostrophe forces Excel's text interpretation, then we
# reset the cell style back to Normal so no stray number-format/style index
# is left behind on the cell.
$ws.Range("E4").Value = "'532832"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = 15.54
$ws.Range("G4").Value = 158.2
$ws.Range("H4").Value = 63215233
